$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A-column date/id values in rows 3-63 (2017xxxx -> 2015xxxx)
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 20000
}

# Reset the view: scroll back to the top and select A1:C63 instead of I64
$ws.Activate()
$ws.Range("A1:C63").Select()
